# 12.8.1 metadata workbook update
# - Update the Indicator text (B4) to drop the parenthetical
#   "(including climate change education)" remark.
# - Leave the active cell selection on B4 (matches the sheet's saved
#   cursor position after the edit).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = "12.8.1 Extent to which (i) global citizenship education and (ii) education for sustainable development are mainstreamed in (a) national education policies; (b) curricula; (c) teacher education; and (d) student assessment"

$ws.Range("B4").Select()
